# Updated cryptos list - refreshed price/volume figures (and fixed the
# WrappedBTC/ShibaInu row ordering) to match the latest GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D/E cells are stored as plain text (prices/percentages are
# formatted strings, not numbers) in the source sheet. Force the cell to
# text ("@") before writing so Excel does not reinterpret e.g. "1.00" or
# "595.46" as a number, then drop back to the Normal style so no stray
# number formatting is left behind on the cell.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "64.808.92"
Set-TextValue "E2" "  +1.72%  "
# Row 3
Set-TextValue "D3" "2.634.65"
Set-TextValue "E3" "  +0.49%  "
# Row 4
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.05%  "
# Row 5
Set-TextValue "D5" "595.46"
Set-TextValue "E5" "  -0.28%  "
# Row 6
Set-TextValue "D6" "154.25"
Set-TextValue "E6" "  +2.13%  "
# Row 7
Set-TextValue "E7" "  -0.01%  "
# Row 8
Set-TextValue "D8" "0.591"
Set-TextValue "E8" "  +0.13%  "
# Row 9
Set-TextValue "E9" "  +5.02%  "
# Row 10
Set-TextValue "E10" "  +2.90%  "
# Row 11
Set-TextValue "D11" "5.79"
Set-TextValue "E11" "  +1.76%  "
# Row 12
Set-TextValue "D12" "0.152"
Set-TextValue "E12" "  +1.30%  "
# Row 13
Set-TextValue "D13" "28.85"
Set-TextValue "E13" "  +3.46%  "
# Row 14
Set-TextValue "D14" "3.113.35"
Set-TextValue "E14" "  +0.54%  "
# Row 15
Set-TextValue "B15" "WrappedBTC"
Set-TextValue "C15" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D15" "64.663.58"
Set-TextValue "E15" "  +1.68%  "
# Row 16
Set-TextValue "B16" "ShibaInu"
Set-TextValue "C16" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D16" "0.0000172"
Set-TextValue "E16" "  +12.83%  "
# Row 17
Set-TextValue "D17" "2.625.69"
Set-TextValue "E17" "  -0.73%  "
# Row 18
Set-TextValue "D18" "12.42"
Set-TextValue "E18" "  +0.64%  "
# Row 19
Set-TextValue "D19" "4.80"
# Row 20
Set-TextValue "D20" "351.95"
Set-TextValue "E20" "  +1.14%  "
# Row 21
Set-TextValue "D21" "7.17"
Set-TextValue "E21" "  +4.16%  "
# Row 22
Set-TextValue "E22" "  +0.22%  "
# Row 23
Set-TextValue "D23" "67.72"
Set-TextValue "E23" "  +1.27%  "
# Row 24
Set-TextValue "D24" "1.71"
Set-TextValue "E24" "  -0.84%  "
# Row 25
Set-TextValue "D25" "9.38"
Set-TextValue "E25" "  +0.72%  "
# Row 26
Set-TextValue "D26" "1.66"
Set-TextValue "E26" "  -1.26%  "
# Row 27
Set-TextValue "D27" "8.28"
Set-TextValue "E27" "  +1.60%  "
# Row 28
Set-TextValue "E28" "  +1.70%  "
# Row 29
Set-TextValue "D29" "541.93"
Set-TextValue "E29" "  -3.20%  "
# Row 30
Set-TextValue "D30" "0.996"
Set-TextValue "E30" "  -0.29%  "
# Row 31
Set-TextValue "D31" "0.0₃0915"
Set-TextValue "E31" "  +7.35%  "
# Row 32
Set-TextValue "D32" "2.06"
Set-TextValue "E32" "  +0.28%  "
# Row 33
Set-TextValue "D33" "1.81"
Set-TextValue "E33" "  +2.89%  "
# Row 34
Set-TextValue "D34" "5.67"
Set-TextValue "E34" "  +7.25%  "
# Row 35
Set-TextValue "D35" "6.25"
Set-TextValue "E35" "  +0.87%  "
# Row 36
Set-TextValue "D36" "0.424"
Set-TextValue "E36" "  +1.83%  "
# Row 37
Set-TextValue "D37" "164.64"
Set-TextValue "E37" "  -2.00%  "
# Row 38
Set-TextValue "D38" "2.02"
Set-TextValue "E38" "  +4.37%  "
# Row 39
Set-TextValue "D39" "20.14"
Set-TextValue "E39" "  +2.84%  "
# Row 40
Set-TextValue "D40" "1.00"
Set-TextValue "E40" "  +0.09%  "
# Row 41
Set-TextValue "E41" "  -0.01%  "
# Row 42
Set-TextValue "D42" "166.62"
Set-TextValue "E42" "  -0.27%  "
# Row 43
Set-TextValue "D43" "41.88"
Set-TextValue "E43" "  +5.53%  "
# Row 44
Set-TextValue "E44" "  +3.75%  "
# Row 45
Set-TextValue "D45" "23.38"
Set-TextValue "E45" "  +6.71%  "
# Row 46
Set-TextValue "D46" "0.0602"
Set-TextValue "E46" "  +1.98%  "
# Row 47
Set-TextValue "E47" "  +10.89%  "
# Row 48
Set-TextValue "D48" "0.643"
Set-TextValue "E48" "  +1.63%  "
# Row 49
Set-TextValue "D49" "0.0251"
Set-TextValue "E49" "  -0.30%  "
# Row 50
Set-TextValue "E50" "  +1.55%  "
# Row 51
Set-TextValue "D51" "19.41"
Set-TextValue "E51" "  -0.05%  "
